$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Type-conversion cells: copy format+value from a cell already bearing the target style/type, then correct numeric cells to the right number ---
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
$ws.Range("I14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("K14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("I14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Plain value updates (same type before/after) ---
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = -9.629629629629
$ws.Range("L16").Value = -25.609756097561
$ws.Range("M16").Value = -14.685314685314
$ws.Range("N16").Value = -77.024482109227
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 85.714285714285
$ws.Range("I17").Value = 219
$ws.Range("J17").Value = 179
$ws.Range("K17").Value = 22.346368715083
$ws.Range("L17").Value = 16.489361702127
$ws.Range("M17").Value = 106.603773584906
$ws.Range("N17").Value = -6.008583690987
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 126
$ws.Range("K18").Value = 15.079365079365
$ws.Range("L18").Value = 23.931623931623
$ws.Range("M18").Value = 1.398601398601
$ws.Range("N18").Value = -85.412474849094
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -52.631578947368
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -37.5
$ws.Range("I19").Value = 643
$ws.Range("J19").Value = 594
$ws.Range("K19").Value = 8.249158249158
$ws.Range("L19").Value = 6.105610561056
$ws.Range("M19").Value = 111.513157894737
$ws.Range("N19").Value = 61.964735516372
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -41.666666666666
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = 3.125
$ws.Range("I20").Value = 296
$ws.Range("J20").Value = 305
$ws.Range("K20").Value = -2.950819672131
$ws.Range("L20").Value = 23.849372384937
$ws.Range("M20").Value = 174.074074074074
$ws.Range("N20").Value = -81.230183893468
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -11.278195488721
$ws.Range("I21").Value = 1437
$ws.Range("J21").Value = 1352
$ws.Range("K21").Value = 6.286982248520
$ws.Range("L21").Value = 7.640449438202
$ws.Range("M21").Value = 76.102941176470
$ws.Range("N21").Value = -61.802232854864
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -50
$ws.Range("M22").Value = -68.75
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 7
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 58
$ws.Range("J23").Value = 72
$ws.Range("K23").Value = -19.444444444444
$ws.Range("L23").Value = -12.121212121212
$ws.Range("M23").Value = 52.631578947368
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 36.842105263157
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -19.672131147541
$ws.Range("I24").Value = 970
$ws.Range("J24").Value = 994
$ws.Range("K24").Value = -2.414486921529
$ws.Range("L24").Value = -5.916585838991
$ws.Range("M24").Value = 39.769452449567
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = -41.666666666666
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -67.241379310344
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = -11.25
$ws.Range("L25").Value = -34.742647058823
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 104.761904761905
$ws.Range("I26").Value = 318
$ws.Range("J26").Value = 269
$ws.Range("K26").Value = 18.215613382899
$ws.Range("L26").Value = 8.163265306122
$ws.Range("M26").Value = -7.826086956521
$ws.Range("L27").Value = -24
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = 37.5
$ws.Range("G29").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = -33.333333333333
